$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").Value = "寻址方式（可选）"
$ws.Range("O2").Value = "静态"
$ws.Range("O3").Value = "动态"
$ws.Range("O3").Select()
